$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was inserted before the existing row 95,
# pushing every subsequent record (old rows 95..192) down by one row
# (new rows 96..193). Insert a new row at 95 to perform that shift.
$ws.Rows.Item(95).Insert()

# Populate the newly inserted row 95 with the new record's data.
$ws.Range("A95").Value = 3
$ws.Range("B95").Value = "Femacal de La Calera"
$ws.Range("C95").Value = "Coquimbo"
$ws.Range("D95").Value = 44539
$ws.Range("E95").Value = 5
$ws.Range("F95").Value = 100112001
$ws.Range("G95").Value = "Berenjena"
$ws.Range("H95").Value = "Sin especificar"
$ws.Range("I95").Value = "Primera"
$ws.Range("J95").Value = 105
$ws.Range("K95").Value = 7500
$ws.Range("L95").Value = 8000
$ws.Range("M95").Value = 7738
$ws.Range("N95").Value = "$/caja 60 unidades"
$ws.Range("O95").Value = "Región de Arica y Parinacota"
$ws.Range("P95").Value = 129
$ws.Range("Q95").Value = 60
$ws.Range("R95").Value = "Hortaliza"
